$wb = $excel.ActiveWorkbook

$wsData      = $wb.Worksheets.Item(1)  # "Data"
$wsDataFinal = $wb.Worksheets.Item(2)  # "Data_final"

# --- Sheet "Data": append a new "Brennstoff allgemein" (general fuel) row ---
# Copy formatting from the row above (row 11) so the new row matches the
# existing style (font/format) used throughout the table.
$wsData.Range("A11").Copy()
$wsData.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

$wsData.Cells.Item(12, 1).Value = "Brennstoff allgemein"
$wsData.Cells.Item(12, 3).Value = 0.9
$wsData.Cells.Item(12, 4).Value = 0
$wsData.Cells.Item(12, 5).Value = "own assumption"

# --- Sheet "Data_final": append the matching new row ---
$wsDataFinal.Range("A11").Copy()
$wsDataFinal.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

$wsDataFinal.Cells.Item(12, 1).Value = "Brennstoff allgemein"
$wsDataFinal.Cells.Item(12, 3).Value = 1
$wsDataFinal.Cells.Item(12, 4).Value = 0

# --- Update selections / active sheet to reflect where the author ended up ---
$wsData.Range("A12:E12").Select()

$wsDataFinal.Activate()
$wsDataFinal.Range("C13").Select()
